# Auto-generated Excel COM-interop script applying the Masamune_Profits
# pricing-refresh diff (scheduled-runner update) across all item sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 12406.1
$ws.Range("J17").Value = 12406.1
$ws.Range("L17").Value = 37218.3
$ws.Range("N17").Value = -37554.3

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 173.30435
$ws.Range("I33").Value = 113.85714
$ws.Range("J33").Value = 265.77777
$ws.Range("K33").Value = 113.85714
$ws.Range("L33").Value = 265.77777
$ws.Range("M33").Value = 115.14286
$ws.Range("N33").Value = -723.7777699999999

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 3650.125
$ws.Range("I51").Value = 2340.2
$ws.Range("J51").Value = 5833.3335
$ws.Range("K51").Value = 2340.2
$ws.Range("L51").Value = 5833.3335
$ws.Range("M51").Value = -1856.2
$ws.Range("N51").Value = -6801.3335

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 3356.25
$ws.Range("I74").Value = 3403.6155
$ws.Range("J74").Value = 3268.2856
$ws.Range("K74").Value = 3403.6155
$ws.Range("L74").Value = 3268.2856
$ws.Range("M74").Value = -2467.6155
$ws.Range("N74").Value = -5140.2856

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 3356.25
$ws.Range("I77").Value = 3403.6155
$ws.Range("J77").Value = 3268.2856
$ws.Range("K77").Value = 17018.0775
$ws.Range("L77").Value = 16341.428
$ws.Range("M77").Value = -12338.0775
$ws.Range("N77").Value = -25701.428

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 9220.700000000001
$ws.Range("I116").Value = 7719
$ws.Range("K116").Value = 7719
$ws.Range("M116").Value = -4277

# Row 128: Nearly There / Kumbhiraskin Grimoire
$ws.Range("H128").Value = 54968.5
$ws.Range("J128").Value = 54968.5
$ws.Range("L128").Value = 54968.5
$ws.Range("N128").Value = -64928.5

# Row 130: Technically Still Magic / Ophiotauroskin Magitek Codex
$ws.Range("H130").Value = 44597.6
$ws.Range("J130").Value = 44597.6
$ws.Range("L130").Value = 44597.6
$ws.Range("N130").Value = -54637.6

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 17858130
$ws.Range("I135").Value = 1022.7778
$ws.Range("K135").Value = 9205.0002
$ws.Range("M135").Value = -6670.0002

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2865.3489
$ws.Range("I137").Value = 1079.2858
$ws.Range("J137").Value = 3442.3845
$ws.Range("K137").Value = 3237.8574
$ws.Range("L137").Value = 10327.1535
$ws.Range("M137").Value = -687.8574000000003
$ws.Range("N137").Value = -15427.1535

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2568.027
$ws.Range("I61").Value = 2375.5625
$ws.Range("J61").Value = 3799.8
$ws.Range("K61").Value = 2375.5625
$ws.Range("L61").Value = 3799.8
$ws.Range("M61").Value = -2163.5625
$ws.Range("N61").Value = -4223.8

# Row 69: The Cut Alembical Cord / Mythrite Alembic
$ws.Range("H69").Value = 99750
$ws.Range("J69").Value = 99750
$ws.Range("L69").Value = 99750
$ws.Range("N69").Value = -101248

# Row 72: Sheer Distill Power (L) / Mythrite Alembic
$ws.Range("H72").Value = 99750
$ws.Range("J72").Value = 99750
$ws.Range("L72").Value = 299250
$ws.Range("N72").Value = -306738

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2036.45
$ws.Range("I74").Value = 1350
$ws.Range("J74").Value = 2875.4443
$ws.Range("K74").Value = 1350
$ws.Range("L74").Value = 2875.4443
$ws.Range("M74").Value = -476
$ws.Range("N74").Value = -4623.4443

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2036.45
$ws.Range("I77").Value = 1350
$ws.Range("J77").Value = 2875.4443
$ws.Range("K77").Value = 6750
$ws.Range("L77").Value = 14377.2215
$ws.Range("M77").Value = -2382
$ws.Range("N77").Value = -23113.2215

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2568.027
$ws.Range("I136").Value = 2375.5625
$ws.Range("J136").Value = 3799.8
$ws.Range("K136").Value = 7126.6875
$ws.Range("L136").Value = 11399.4
$ws.Range("M136").Value = -4576.6875
$ws.Range("N136").Value = -16499.4

# Row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting
$ws.Range("H138").Value = 44970
$ws.Range("J138").Value = 44970
$ws.Range("L138").Value = 44970
$ws.Range("N138").Value = -55250

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2489.48
$ws.Range("I105").Value = 2481.9333
$ws.Range("J105").Value = 2500.8
$ws.Range("K105").Value = 2481.9333
$ws.Range("L105").Value = 2500.8
$ws.Range("M105").Value = -734.9333000000001
$ws.Range("N105").Value = -5994.8

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 2686.111
$ws.Range("I62").Value = 2596
$ws.Range("J62").Value = 2798.75
$ws.Range("K62").Value = 2596
$ws.Range("L62").Value = 2798.75
$ws.Range("M62").Value = -1972
$ws.Range("N62").Value = -4046.75

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 2686.111
$ws.Range("I65").Value = 2596
$ws.Range("J65").Value = 2798.75
$ws.Range("K65").Value = 12980
$ws.Range("L65").Value = 13993.75
$ws.Range("M65").Value = -9860
$ws.Range("N65").Value = -20233.75

# Row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws.Range("H80").Value = 31622.8
$ws.Range("J80").Value = 31622.8
$ws.Range("L80").Value = 31622.8
$ws.Range("N80").Value = -33868.8

# Row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws.Range("H83").Value = 31622.8
$ws.Range("J83").Value = 31622.8
$ws.Range("L83").Value = 94868.39999999999
$ws.Range("N83").Value = -106100.4

# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 1820.3334
$ws.Range("I86").Value = 1801.4546
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1801.4546
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -678.4546
$ws.Range("N86").Value = -4096

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 1820.3334
$ws.Range("I89").Value = 1801.4546
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 9007.273000000001
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -3391.273000000001
$ws.Range("N89").Value = -20482

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy / Birch Syrup
$ws.Range("H86").Value = 1579.909
$ws.Range("I86").Value = 693.3333
$ws.Range("J86").Value = 1912.375
$ws.Range("K86").Value = 2079.9999
$ws.Range("L86").Value = 5737.125
$ws.Range("M86").Value = -893.9998999999998
$ws.Range("N86").Value = -8109.125

# Row 89: Luxury Spillover (L) / Birch Syrup
$ws.Range("H89").Value = 1579.909
$ws.Range("I89").Value = 693.3333
$ws.Range("J89").Value = 1912.375
$ws.Range("K89").Value = 6239.9997
$ws.Range("L89").Value = 17211.375
$ws.Range("M89").Value = -311.9997000000003
$ws.Range("N89").Value = -29067.375

$ws = $wb.Worksheets.Item("GSM")
# Row 15: The Tusk at Hand / Fang Earrings
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 5013.5137
$ws.Range("I70").Value = 4988
$ws.Range("K70").Value = 4988
$ws.Range("M70").Value = -4718

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 5013.5137
$ws.Range("I73").Value = 4988
$ws.Range("K73").Value = 4988
$ws.Range("M73").Value = -4052

# Row 81: The Grander Temple / Dragon Fang Earrings
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84: Man with a Dragon Earring (L) / Dragon Fang Earrings
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 7175.087
$ws.Range("I97").Value = 1615.5
$ws.Range("K97").Value = 1615.5
$ws.Range("M97").Value = -1119.5

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2867.3333
$ws.Range("I102").Value = 2835.2942
$ws.Range("K102").Value = 2835.2942
$ws.Range("M102").Value = -1213.2942

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1671.0625
$ws.Range("I122").Value = 1673.9166
$ws.Range("J122").Value = 1662.5
$ws.Range("K122").Value = 5021.7498
$ws.Range("L122").Value = 4987.5
$ws.Range("M122").Value = -2571.7498
$ws.Range("N122").Value = -9887.5

# Row 128: To Fight at Her Side / Manganese Rapier
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 63: From Mud to Mourning / Archaeoskin Jackboots of Gathering
$ws.Range("H63").Value = 11500
$ws.Range("J63").Value = 11500
$ws.Range("L63").Value = 11500
$ws.Range("N63").Value = -12998

# Row 66: These Boots Are Made for Hawkin' (L) / Archaeoskin Jackboots of Gathering
$ws.Range("H66").Value = 11500
$ws.Range("J66").Value = 11500
$ws.Range("L66").Value = 34500
$ws.Range("N66").Value = -41988

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3451.1458
$ws.Range("I132").Value = 3301.7334
$ws.Range("J132").Value = 3700.1667
$ws.Range("K132").Value = 9905.200199999999
$ws.Range("L132").Value = 11100.5001
$ws.Range("M132").Value = -7375.200199999999
$ws.Range("N132").Value = -16160.5001

# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 25560.715
$ws.Range("J133").Value = 25560.715
$ws.Range("L133").Value = 25560.715
$ws.Range("N133").Value = -30620.715

# Row 134: Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 48482
$ws.Range("J134").Value = 48482
$ws.Range("L134").Value = 48482
$ws.Range("N134").Value = -58622

# Row 138: Freezing Toes / Gomphotherium Boots of Striking
$ws.Range("H138").Value = 54784.4
$ws.Range("J138").Value = 54784.4
$ws.Range("L138").Value = 54784.4
$ws.Range("N138").Value = -65064.4

$ws = $wb.Worksheets.Item("WVR")
# Row 86: Felt for the Fallen / Chimerical Felt
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89: Blinded Veil of Vigilance (L) / Chimerical Felt
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1331.3334
$ws.Range("I132").Value = 949.62964
$ws.Range("J132").Value = 3049
$ws.Range("K132").Value = 2848.88892
$ws.Range("L132").Value = 9147
$ws.Range("M132").Value = -318.8889199999999
$ws.Range("N132").Value = -14207
